$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Weekly price-list refresh: row 140 receives the newest weeks figures,
# rows 141-160 shift down to what used to be in the row above them, and
# row 161 is newly appended holding the data that used to live in row 160.

$ws.Range("D140").Value = 44951
$ws.Range("J140").Value = 15
$ws.Range("K140").Value = 24000
$ws.Range("L140").Value = 24000
$ws.Range("M140").Value = 24000
$ws.Range("N140").Value = "$/malla 20 kilos"
$ws.Range("O140").Value = "Perú"
$ws.Range("P140").Value = 1200
$ws.Range("Q140").Value = 20

$ws.Range("D141").Value = 44663
$ws.Range("J141").Value = 30
$ws.Range("K141").Value = 18000
$ws.Range("L141").Value = 18000
$ws.Range("M141").Value = 18000
$ws.Range("N141").Value = "$/malla 20 kilos"
$ws.Range("O141").Value = "Perú"
$ws.Range("P141").Value = 900
$ws.Range("Q141").Value = 20

$ws.Range("D142").Value = 44690
$ws.Range("J142").Value = 15
$ws.Range("K142").Value = 18000
$ws.Range("L142").Value = 18000
$ws.Range("M142").Value = 18000
$ws.Range("N142").Value = "$/malla 20 kilos"
$ws.Range("O142").Value = "Perú"
$ws.Range("P142").Value = 900
$ws.Range("Q142").Value = 20

$ws.Range("D143").Value = 44607
$ws.Range("J143").Value = 40
$ws.Range("K143").Value = 18000
$ws.Range("L143").Value = 18000
$ws.Range("M143").Value = 18000
$ws.Range("N143").Value = "$/malla 20 kilos"
$ws.Range("O143").Value = "Perú"
$ws.Range("P143").Value = 900
$ws.Range("Q143").Value = 20

$ws.Range("D144").Value = 44529
$ws.Range("J144").Value = 15
$ws.Range("K144").Value = 20000
$ws.Range("L144").Value = 20000
$ws.Range("M144").Value = 20000
$ws.Range("N144").Value = "$/malla 20 kilos"
$ws.Range("O144").Value = "Perú"
$ws.Range("P144").Value = 1000
$ws.Range("Q144").Value = 20

$ws.Range("D145").Value = 44901
$ws.Range("J145").Value = 15
$ws.Range("K145").Value = 24000
$ws.Range("L145").Value = 24000
$ws.Range("M145").Value = 24000
$ws.Range("N145").Value = "$/malla 20 kilos"
$ws.Range("O145").Value = "Perú"
$ws.Range("P145").Value = 1200
$ws.Range("Q145").Value = 20

$ws.Range("D146").Value = 44637
$ws.Range("J146").Value = 30
$ws.Range("K146").Value = 18000
$ws.Range("L146").Value = 18000
$ws.Range("M146").Value = 18000
$ws.Range("N146").Value = "$/malla 20 kilos"
$ws.Range("O146").Value = "Perú"
$ws.Range("P146").Value = 900
$ws.Range("Q146").Value = 20

$ws.Range("D147").Value = 44855
$ws.Range("J147").Value = 40
$ws.Range("K147").Value = 25000
$ws.Range("L147").Value = 25000
$ws.Range("M147").Value = 25000
$ws.Range("N147").Value = "$/malla 20 kilos"
$ws.Range("O147").Value = "Perú"
$ws.Range("P147").Value = 1250
$ws.Range("Q147").Value = 20

$ws.Range("D148").Value = 44616
$ws.Range("J148").Value = 80
$ws.Range("K148").Value = 15000
$ws.Range("L148").Value = 15000
$ws.Range("M148").Value = 15000
$ws.Range("N148").Value = "$/malla 20 kilos"
$ws.Range("O148").Value = "Perú"
$ws.Range("P148").Value = 750
$ws.Range("Q148").Value = 20

$ws.Range("D149").Value = 44650
$ws.Range("J149").Value = 20
$ws.Range("K149").Value = 18000
$ws.Range("L149").Value = 18000
$ws.Range("M149").Value = 18000
$ws.Range("N149").Value = "$/malla 20 kilos"
$ws.Range("O149").Value = "Perú"
$ws.Range("P149").Value = 900
$ws.Range("Q149").Value = 20

$ws.Range("D150").Value = 44340
$ws.Range("J150").Value = 40
$ws.Range("K150").Value = 18000
$ws.Range("L150").Value = 18000
$ws.Range("M150").Value = 18000
$ws.Range("N150").Value = "$/malla 20 kilos"
$ws.Range("O150").Value = "Perú"
$ws.Range("P150").Value = 900
$ws.Range("Q150").Value = 20

$ws.Range("D151").Value = 44908
$ws.Range("J151").Value = 20
$ws.Range("K151").Value = 24000
$ws.Range("L151").Value = 24000
$ws.Range("M151").Value = 24000
$ws.Range("N151").Value = "$/malla 20 kilos"
$ws.Range("O151").Value = "Perú"
$ws.Range("P151").Value = 1200
$ws.Range("Q151").Value = 20

$ws.Range("D152").Value = 44826
$ws.Range("J152").Value = 100
$ws.Range("K152").Value = 20000
$ws.Range("L152").Value = 20000
$ws.Range("M152").Value = 20000
$ws.Range("N152").Value = "$/malla 20 kilos"
$ws.Range("O152").Value = "Perú"
$ws.Range("P152").Value = 1000
$ws.Range("Q152").Value = 20

$ws.Range("D153").Value = 44757
$ws.Range("J153").Value = 40
$ws.Range("K153").Value = 20000
$ws.Range("L153").Value = 20000
$ws.Range("M153").Value = 20000
$ws.Range("N153").Value = "$/malla 20 kilos"
$ws.Range("O153").Value = "Perú"
$ws.Range("P153").Value = 1000
$ws.Range("Q153").Value = 20

$ws.Range("D154").Value = 44838
$ws.Range("J154").Value = 30
$ws.Range("K154").Value = 20000
$ws.Range("L154").Value = 20000
$ws.Range("M154").Value = 20000
$ws.Range("N154").Value = "$/malla 20 kilos"
$ws.Range("O154").Value = "Perú"
$ws.Range("P154").Value = 1000
$ws.Range("Q154").Value = 20

$ws.Range("D155").Value = 44671
$ws.Range("J155").Value = 25
$ws.Range("K155").Value = 18000
$ws.Range("L155").Value = 18000
$ws.Range("M155").Value = 18000
$ws.Range("N155").Value = "$/malla 20 kilos"
$ws.Range("O155").Value = "Perú"
$ws.Range("P155").Value = 900
$ws.Range("Q155").Value = 20

$ws.Range("D156").Value = 44657
$ws.Range("J156").Value = 15
$ws.Range("K156").Value = 18000
$ws.Range("L156").Value = 18000
$ws.Range("M156").Value = 18000
$ws.Range("N156").Value = "$/caja 15 kilos granel"
$ws.Range("O156").Value = "Perú"
$ws.Range("P156").Value = 1200
$ws.Range("Q156").Value = 15

$ws.Range("D157").Value = 44636
$ws.Range("J157").Value = 50
$ws.Range("K157").Value = 18000
$ws.Range("L157").Value = 18000
$ws.Range("M157").Value = 18000
$ws.Range("N157").Value = "$/malla 20 kilos"
$ws.Range("O157").Value = "Perú"
$ws.Range("P157").Value = 900
$ws.Range("Q157").Value = 20

$ws.Range("D158").Value = 44810
$ws.Range("J158").Value = 30
$ws.Range("K158").Value = 20000
$ws.Range("L158").Value = 20000
$ws.Range("M158").Value = 20000
$ws.Range("N158").Value = "$/malla 20 kilos"
$ws.Range("O158").Value = "Perú"
$ws.Range("P158").Value = 1000
$ws.Range("Q158").Value = 20

$ws.Range("D159").Value = 44175
$ws.Range("J159").Value = 20
$ws.Range("K159").Value = 20000
$ws.Range("L159").Value = 20000
$ws.Range("M159").Value = 20000
$ws.Range("N159").Value = "$/caja 15 kilos granel"
$ws.Range("O159").Value = "Región de Arica y Parinacota"
$ws.Range("P159").Value = 1333
$ws.Range("Q159").Value = 15

$ws.Range("D160").Value = 44795
$ws.Range("J160").Value = 50
$ws.Range("K160").Value = 20000
$ws.Range("L160").Value = 20000
$ws.Range("M160").Value = 20000
$ws.Range("N160").Value = "$/malla 20 kilos"
$ws.Range("O160").Value = "Perú"
$ws.Range("P160").Value = 1000
$ws.Range("Q160").Value = 20

$ws.Range("D161").Value = 44595
$ws.Range("J161").Value = 50
$ws.Range("K161").Value = 18000
$ws.Range("L161").Value = 18000
$ws.Range("M161").Value = 18000
$ws.Range("N161").Value = "$/malla 20 kilos"
$ws.Range("O161").Value = "Perú"
$ws.Range("P161").Value = 900
$ws.Range("Q161").Value = 20

# Row 161 is brand-new: fill in the columns that stay constant across
# every record for this sheet (same values already used by every other
# row, e.g. row 160).
$ws.Range("A161").Value = 10
$ws.Range("B161").Value = "Vega Modelo de Temuco"
$ws.Range("C161").Value = "La Araucanía"
$ws.Range("E161").Value = 9
$ws.Range("F161").Value = 100114002
$ws.Range("G161").Value = "Camote"
$ws.Range("H161").Value = "Sin especificar"
$ws.Range("I161").Value = "Primera"
$ws.Range("R161").Value = "Hortaliza"

# Match the date/time display style already used by column D.
$ws.Range("D161").NumberFormat = "YYYY-MM-DD HH:MM:SS"
